$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reference cells with known good styles for columns B (style 2) and C (style 3)
$styleRefB = $ws.Range("B2")
$styleRefC = $ws.Range("C2")

# Row 10
$ws.Range("A10").Value = 'Objetivos:'
$ws.Range("B10").Value = '519033 - Carlos Yujiro Shigue'
$ws.Range("C10").Value = '519033 - Carlos Yujiro Shigue'
$ws.Rows.Item(10).RowHeight = 60

# Row 11
$ws.Range("A11").Value = 'Objectives:'
$ws.Rows.Item(11).RowHeight = 60

# Row 12
$ws.Range("A12").Value = 'Docentes responsáveis:'

# Row 13
$ws.Range("A13").Value = 'Programa resumido:'
$ws.Range("B13").Value = '519033 - Carlos Yujiro Shigue'
$ws.Range("C13").Value = '519033 - Carlos Yujiro Shigue'
$ws.Rows.Item(13).RowHeight = 60

# Row 14
$ws.Range("B14").Clear()
$ws.Range("C14").Clear()
$ws.Range("A14").Value = 'Short syllabus:'
$ws.Rows.Item(14).RowHeight = 60

# Row 15
$ws.Range("A15").Value = 'Programa:'
$ws.Range("B15").Value = '1033242 - Fábio Herbst Florenzano'
$ws.Range("C15").Value = '1033242 - Fábio Herbst Florenzano'
$ws.Rows.Item(15).RowHeight = 120

# Row 16
$ws.Range("B16").Clear()
$ws.Range("C16").Clear()
$ws.Range("A16").Value = 'Syllabus:'
$ws.Rows.Item(16).RowHeight = 120

# Row 17
$ws.Range("A17").Value = 'Avaliação:'
$ws.Rows.Item(17).AutoFit()

# Row 18
$ws.Range("B18").Clear()
$ws.Range("C18").Clear()
$ws.Range("A18").Value = 'Método:'
$ws.Rows.Item(18).RowHeight = 60

# Row 19
$ws.Range("A19").Value = 'Critério:'
$ws.Range("B19").Value = 'De acordo com a atual ementa da disciplina propõe-se o uso de uma nova metodologia de ensino com o intuito de abordar o conteúdo de forma mais prática e contextualizada para que o aluno consiga relacionar os conhecimentos teóricos vistos em sala de aula com as outras disciplinas do curso. Assim, avaliação do aluno será feita através de uma prova escrita e por uma apresentação final com base nas atividades práticas desenvolvidas.'
$ws.Range("C19").Value = 'De acordo com a atual ementa da disciplina propõe-se o uso de uma nova metodologia de ensino com o intuito de abordar o conteúdo de forma mais prática e contextualizada para que o aluno consiga relacionar os conhecimentos teóricos vistos em sala de aula com as outras disciplinas do curso. Assim, avaliação do aluno será feita através de uma prova escrita e por uma apresentação final com base nas atividades práticas desenvolvidas.'
$styleRefB.Copy()
$ws.Range("B19").PasteSpecial(-4122)
$styleRefC.Copy()
$ws.Range("C19").PasteSpecial(-4122)
$ws.Rows.Item(19).RowHeight = 60

# Row 20
$ws.Range("A20").Value = 'Norma de recuperação:'
$ws.Range("B20").Value = 'A nota final será calculada como descrita a seguir: NF= (0,4*Avaliação escrita + 0,6 *Apresentação final)'
$ws.Range("C20").Value = 'A nota final será calculada como descrita a seguir: NF= (0,4*Avaliação escrita + 0,6 *Apresentação final)'
$styleRefB.Copy()
$ws.Range("B20").PasteSpecial(-4122)
$styleRefC.Copy()
$ws.Range("C20").PasteSpecial(-4122)
$ws.Rows.Item(20).RowHeight = 60

# Row 21
$ws.Range("A21").Value = 'Bibliografia:'
$ws.Range("B21").Value = 'Devido a cunho prático da disciplina não haverá recuperação.'
$ws.Range("C21").Value = 'Devido a cunho prático da disciplina não haverá recuperação.'
$ws.Rows.Item(21).RowHeight = 120

# Row 22
$ws.Range("B22").Clear()
$ws.Range("C22").Clear()
$ws.Range("A22").Value = 'Requisitos:'
$ws.Rows.Item(22).AutoFit()

# Row 23
$ws.Range("A23").Clear()
$ws.Range("B23").Value = 'LOM3011 -  Ensaios Mecânicos  (Requisito fraco)'
$ws.Range("C23").Value = 'LOM3011 -  Ensaios Mecânicos  (Requisito fraco)'
$ws.Rows.Item(23).RowHeight = 30

# Row 24
$ws.Range("A24").Clear()
$ws.Range("B24").Value = 'LOM3046 -  Técnicas de Análise Microestrutural  (Requisito fraco)'
$ws.Range("C24").Value = 'LOM3046 -  Técnicas de Análise Microestrutural  (Requisito fraco)'
$ws.Rows.Item(24).RowHeight = 30

$excel.CutCopyMode = $false

# Remove trailing rows 25-27 entirely (content consolidated into rows above)
$ws.Rows.Item(25).Delete()
$ws.Rows.Item(25).Delete()
$ws.Rows.Item(25).Delete()
